$wb = $excel.ActiveWorkbook

# --- SUPPLIER CONFIG sheet ---
$ws1 = $wb.Worksheets.Item("SUPPLIER CONFIG")
$ws1.Range("D6").Value = 125
$ws1.Range("D7").Value = 100
$ws1.Range("D8").Value = 140
$ws1.Range("D9").Value = 330
$ws1.Range("D10").Value = 264
$ws1.Range("D11").Value = 370

$ws1.Range("B16").Value = 60
$ws1.Range("C16").Value = 1
$ws1.Range("B17").Value = 7
$ws1.Range("C17").Value = 100
$ws1.Range("B18").Value = 36
$ws1.Range("C18").Value = 30
$ws1.Range("B19").Value = 24
$ws1.Range("C19").Value = 60
$ws1.Range("B20").Value = 30
$ws1.Range("C20").Value = 100
$ws1.Range("B21").Value = 28
$ws1.Range("C21").Value = 150

# --- COST ANALYSIS sheet ---
$ws2 = $wb.Worksheets.Item("COST ANALYSIS")
$ws2.Range("B24").Formula = "=AVERAGE('SUPPLIER CONFIG'!D6:D8)"
$ws2.Range("C24").Value = 125
$ws2.Range("B25").Formula = "=AVERAGE('SUPPLIER CONFIG'!D9:D11)"
$ws2.Range("C25").Value = 330

# --- CASH FLOW PREVIEW sheet ---
$ws4 = $wb.Worksheets.Item("CASH FLOW PREVIEW")
$ws4.Range("B6").Formula = "='MRP ENGINE'!B18+'MRP ENGINE'!B30"
$ws4.Range("C6").Formula = "='MRP ENGINE'!C18+'MRP ENGINE'!C30"
$ws4.Range("D6").Formula = "='MRP ENGINE'!D18+'MRP ENGINE'!D30"
$ws4.Range("E6").Formula = "='MRP ENGINE'!E18+'MRP ENGINE'!E30"
$ws4.Range("F6").Formula = "='MRP ENGINE'!F18+'MRP ENGINE'!F30"
$ws4.Range("G6").Formula = "='MRP ENGINE'!G18+'MRP ENGINE'!G30"
$ws4.Range("H6").Formula = "='MRP ENGINE'!H18+'MRP ENGINE'!H30"
$ws4.Range("I6").Formula = "='MRP ENGINE'!I18+'MRP ENGINE'!I30"

$ws4.Range("B7").Formula = "='MRP ENGINE'!B19+'MRP ENGINE'!B31"
$ws4.Range("C7").Formula = "='MRP ENGINE'!C19+'MRP ENGINE'!C31"
$ws4.Range("D7").Formula = "='MRP ENGINE'!D19+'MRP ENGINE'!D31"
$ws4.Range("E7").Formula = "='MRP ENGINE'!E19+'MRP ENGINE'!E31"
$ws4.Range("F7").Formula = "='MRP ENGINE'!F19+'MRP ENGINE'!F31"
$ws4.Range("G7").Formula = "='MRP ENGINE'!G19+'MRP ENGINE'!G31"
$ws4.Range("H7").Formula = "='MRP ENGINE'!H19+'MRP ENGINE'!H31"
$ws4.Range("I7").Formula = "='MRP ENGINE'!I19+'MRP ENGINE'!I31"

$ws4.Range("B8").Formula = "='MRP ENGINE'!B20+'MRP ENGINE'!B32"
$ws4.Range("C8").Formula = "='MRP ENGINE'!C20+'MRP ENGINE'!C32"
$ws4.Range("D8").Formula = "='MRP ENGINE'!D20+'MRP ENGINE'!D32"
$ws4.Range("E8").Formula = "='MRP ENGINE'!E20+'MRP ENGINE'!E32"
$ws4.Range("F8").Formula = "='MRP ENGINE'!F20+'MRP ENGINE'!F32"
$ws4.Range("G8").Formula = "='MRP ENGINE'!G20+'MRP ENGINE'!G32"
$ws4.Range("H8").Formula = "='MRP ENGINE'!H20+'MRP ENGINE'!H32"
$ws4.Range("I8").Formula = "='MRP ENGINE'!I20+'MRP ENGINE'!I32"

# --- UPLOAD READY PROCUREMENT sheet ---
$ws5 = $wb.Worksheets.Item("UPLOAD READY PROCUREMENT")
$ws5.Range("D6").Formula = "='MRP ENGINE'!B18"
$ws5.Range("E6").Formula = "='MRP ENGINE'!C18"
$ws5.Range("F6").Formula = "='MRP ENGINE'!D18"
$ws5.Range("G6").Formula = "='MRP ENGINE'!E18"
$ws5.Range("H6").Formula = "='MRP ENGINE'!F18"
$ws5.Range("I6").Formula = "='MRP ENGINE'!G18"
$ws5.Range("J6").Formula = "='MRP ENGINE'!H18"
$ws5.Range("K6").Formula = "='MRP ENGINE'!I18"

$ws5.Range("D7").Formula = "='MRP ENGINE'!B30"
$ws5.Range("E7").Formula = "='MRP ENGINE'!C30"
$ws5.Range("F7").Formula = "='MRP ENGINE'!D30"
$ws5.Range("G7").Formula = "='MRP ENGINE'!E30"
$ws5.Range("H7").Formula = "='MRP ENGINE'!F30"
$ws5.Range("I7").Formula = "='MRP ENGINE'!G30"
$ws5.Range("J7").Formula = "='MRP ENGINE'!H30"
$ws5.Range("K7").Formula = "='MRP ENGINE'!I30"

$ws5.Range("D8").Formula = "='MRP ENGINE'!B19"
$ws5.Range("E8").Formula = "='MRP ENGINE'!C19"
$ws5.Range("F8").Formula = "='MRP ENGINE'!D19"
$ws5.Range("G8").Formula = "='MRP ENGINE'!E19"
$ws5.Range("H8").Formula = "='MRP ENGINE'!F19"
$ws5.Range("I8").Formula = "='MRP ENGINE'!G19"
$ws5.Range("J8").Formula = "='MRP ENGINE'!H19"
$ws5.Range("K8").Formula = "='MRP ENGINE'!I19"

$ws5.Range("D9").Formula = "='MRP ENGINE'!B31"
$ws5.Range("E9").Formula = "='MRP ENGINE'!C31"
$ws5.Range("F9").Formula = "='MRP ENGINE'!D31"
$ws5.Range("G9").Formula = "='MRP ENGINE'!E31"
$ws5.Range("H9").Formula = "='MRP ENGINE'!F31"
$ws5.Range("I9").Formula = "='MRP ENGINE'!G31"
$ws5.Range("J9").Formula = "='MRP ENGINE'!H31"
$ws5.Range("K9").Formula = "='MRP ENGINE'!I31"

$ws5.Range("D10").Formula = "='MRP ENGINE'!B20"
$ws5.Range("E10").Formula = "='MRP ENGINE'!C20"
$ws5.Range("F10").Formula = "='MRP ENGINE'!D20"
$ws5.Range("G10").Formula = "='MRP ENGINE'!E20"
$ws5.Range("H10").Formula = "='MRP ENGINE'!F20"
$ws5.Range("I10").Formula = "='MRP ENGINE'!G20"
$ws5.Range("J10").Formula = "='MRP ENGINE'!H20"
$ws5.Range("K10").Formula = "='MRP ENGINE'!I20"

$ws5.Range("D11").Formula = "='MRP ENGINE'!B32"
$ws5.Range("E11").Formula = "='MRP ENGINE'!C32"
$ws5.Range("F11").Formula = "='MRP ENGINE'!D32"
$ws5.Range("G11").Formula = "='MRP ENGINE'!E32"
$ws5.Range("H11").Formula = "='MRP ENGINE'!F32"
$ws5.Range("I11").Formula = "='MRP ENGINE'!G32"
$ws5.Range("J11").Formula = "='MRP ENGINE'!H32"
$ws5.Range("K11").Formula = "='MRP ENGINE'!I32"

# --- CROSS REFERENCE sheet ---
$ws6 = $wb.Worksheets.Item("CROSS REFERENCE")
$ws6.Range("B6").Value = 0
$ws6.Range("B11").Value = 0
